$d = $word.ActiveDocument

# --- Change 1: add ", Crystal Reports" after "Agile Development, SDLC" in the skills table ---
$d.Content.Find.Execute(
    "Agile Development, SDLC", $false, $false, $false, $false, $false,
    $true, 1, $false, "Agile Development, SDLC, Crystal Reports", 2)

# --- Change 2: expand the "Built numerous..." bullet paragraph ---
$oldText = [char]0x2022 + " Built numerous Excel direct connect reports and designed Crystal Reports to give actionable insights into customer data such as terminations, new hires, and moves on a weekly, monthly, and yearly basis."
$newText = [char]0x2022 + " Built Excel direct connect reports and designed Crystal Reports using SQL queries against customer data to produce actionable insights into customer data such as rentable square feet, terminations, new hires, and moves on a weekly, monthly, and yearly basis."

$d.Content.Find.Execute(
    $oldText, $false, $false, $false, $false, $false,
    $true, 1, $false, $newText, 2)

# --- Change 3: re-add the _GoBack bookmark on the now-empty paragraph that follows ---
# Locate the paragraph that now contains the rewritten bullet, then bookmark the
# start of the blank paragraph immediately after it.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*yearly basis.*") {
        $target = $d.Paragraphs.Item($i + 1)
        $bmRange = $d.Range($target.Range.Start - 1, $target.Range.End)
        $d.Bookmarks.Add("_GoBack", $bmRange)
        break
    }
}
